# Insert a new row at position 686, shifting all existing rows (686..742)
# down by one (to 687..743). This brings the used range from A1:T742 to
# A1:T743.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("686:686").Insert()

# Populate the newly inserted row 686 with its data. Columns A, B, C,
# E-L, N-T keep the same values that row 686 held before the insert
# (which is now sitting in row 687); only D (Fecha) and M (Volumen)
# take new values.
$ws.Cells.Item(686, 1).Value = 3
$ws.Cells.Item(686, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(686, 3).Value = "Coquimbo"
$ws.Cells.Item(686, 4).Value = 45223
$ws.Cells.Item(686, 5).Value = 5
$ws.Cells.Item(686, 6).Value = "Fruta"
$ws.Cells.Item(686, 7).Value = 100108
$ws.Cells.Item(686, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(686, 9).Value = 100108002
$ws.Cells.Item(686, 10).Value = "Mango"
$ws.Cells.Item(686, 11).Value = "Sin especificar"
$ws.Cells.Item(686, 12).Value = "Primera"
$ws.Cells.Item(686, 13).Value = 220
$ws.Cells.Item(686, 14).Value = 9000
$ws.Cells.Item(686, 15).Value = 9000
$ws.Cells.Item(686, 16).Value = 9000
$ws.Cells.Item(686, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(686, 18).Value = "Brasil"
$ws.Cells.Item(686, 19).Value = 2250
$ws.Cells.Item(686, 20).Value = 4
